$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.194.47'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.891.54'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '''320.71'
$ws.Range("E5").Value = '  -2.76%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").Value = '''0.5056'
$ws.Range("E7").Value = '  -3.26%  '
$ws.Range("D8").Value = '''0.4025'
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("D9").Value = '''0.08280'
$ws.Range("E9").Value = '  -2.81%  '
$ws.Range("D10").Value = '''1.107'
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("D11").Value = '''42.20'
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").Value = '''24.16'
$ws.Range("E12").Value = '  +7.80%  '
$ws.Range("D13").Value = '''6.379'
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").Value = '1.884.51'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").Value = '''7.310'
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = '''1.004'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").Value = '''92.65'
$ws.Range("E17").Value = '  -2.50%  '
$ws.Range("D18").Value = '''0.00001097'
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").Value = '''0.06457'
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").Value = '''18.40'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").Value = '''5.913'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("D23").Value = '30.216.10'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '''11.26'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").Value = '''2.189'
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("D26").Value = '2.108.06'
$ws.Range("E26").Value = '  -1.26%  '
$ws.Range("D27").Value = '''21.57'
$ws.Range("E27").Value = '  +2.23%  '
$ws.Range("D28").Value = '''160.71'
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Value = '''2.275'
$ws.Range("E29").Value = '  -5.85%  '
$ws.Range("D30").Value = '''128.76'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").Value = '''1.110'
$ws.Range("E31").Value = '  +2.74%  '
$ws.Range("D32").Value = '''0.1042'
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("D33").Value = '''5.998'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  +2.84%  '
$ws.Range("D35").Value = '''0.02446'
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").Value = '''5.322'
$ws.Range("E36").Value = '  +2.85%  '
$ws.Range("D37").Value = '''0.06429'
$ws.Range("E37").Value = '  -2.28%  '
$ws.Range("D38").Value = '''0.2155'
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").Value = '''1.186'
$ws.Range("E39").Value = '  -3.47%  '
$ws.Range("D40").Value = '''8.600'
$ws.Range("E40").Value = '  -3.00%  '
$ws.Range("D41").Value = '''0.6381'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").Value = '''11.37'
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("D43").Value = '''1.213'
$ws.Range("D44").Value = '''1.001'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '''13.26'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").Value = '''0.5962'
$ws.Range("E46").Value = '  -2.99%  '
$ws.Range("D47").Value = '''2.147'
$ws.Range("E47").Value = '  +3.25%  '
$ws.Range("D48").Value = '''3.639'
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("D49").Value = '''123.45'
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("D50").Value = '''1.213'
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("D51").Value = '''78.57'
$ws.Range("E51").Value = '  -1.20%  '
